$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 146; this shifts existing rows 146:225 down to 147:226
$ws.Rows("146").Insert()

# Populate the newly inserted row 146 with the new record's data
$ws.Cells.Item(146, 1).Value = 4
$ws.Cells.Item(146, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(146, 3).Value = "Los Lagos"
$ws.Cells.Item(146, 4).Value = 44460
$ws.Cells.Item(146, 5).Value = 10
$ws.Cells.Item(146, 6).Value = 100114001
$ws.Cells.Item(146, 7).Value = "Papa"
$ws.Cells.Item(146, 8).Value = "Rodeo"
$ws.Cells.Item(146, 9).Value = "1a (guarda)"
$ws.Cells.Item(146, 10).Value = 300
$ws.Cells.Item(146, 11).Value = 8000
$ws.Cells.Item(146, 12).Value = 8000
$ws.Cells.Item(146, 13).Value = 8000
$ws.Cells.Item(146, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(146, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(146, 16).Value = 320
$ws.Cells.Item(146, 17).Value = 25
$ws.Cells.Item(146, 18).Value = "Hortaliza"
